# Add a new worksheet "8_" at the end of the workbook (after sheet "7"),
# make it the active sheet/tab, and populate it with the "Correct order of
# definitions" matching question (Area of the Plot / definitions table).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "8_"

# Column widths (approximate characters; engine stores width = chars + 5/6)
$ws.Columns.Item(1).ColumnWidth = 21.45
$ws.Columns.Item(2).ColumnWidth = 27.74
$ws.Columns.Item(3).ColumnWidth = 34.88
$ws.Columns.Item(4).ColumnWidth = 33.59

# Row 1 - header
$ws.Range("A1").Value = "Area of the Plot"
$ws.Range("B1").Value = "Correct order of definitions"
$ws.Range("C1").Value = "Definitions"

# Row 2
$ws.Range("A2").Value = "Population: 0-3 billion"
$ws.Range("B2").Value = "C"
$ws.Range("C2").Value = "In this section, the alpha and beta terms are about the same size (growth is very small and sometimes negative)"

# Row 3
$ws.Range("A3").Value = "Population: 3-7 billion"
$ws.Range("B3").Value = "D"
$ws.Range("C3").Value = "In this section, the beta term is bigger than the alpha term (growth is negative)"

# Row 4
$ws.Range("A4").Value = "Population 7-10 billion"
$ws.Range("B4").Value = "D"
$ws.Range("C4").Value = "This section is dominated by the alpha term (growth is small but always positive)"
$ws.Range("D4").Value = "Don't be fooled by the downward trend of the curve: the population growth here (as shown on the y-axis) is still positive and large"

# Row 5
$ws.Range("A5").Value = "Population: 12-14 billion"
$ws.Range("B5").Value = "A"
$ws.Range("C5").Value = "In this section, the alpha term is significantly bigger than the beta term (growth is large)"
$ws.Range("D5").Value = "A population in this area is nearing equilibrium: no positive or negative growth."

# Row 6
$ws.Range("A6").Value = "Population: Above 14 billion"
$ws.Range("B6").Value = "B"
$ws.Range("D6").Value = "If the population is quite large, then it will actually decline"

# Formatting: column A / C wrap text; column B wrap + centered (both axes);
# D4 is wrapped but D5/D6 are left with the default (no) style.
$ws.Range("A1:A6").WrapText = $true
$ws.Range("C1:C5").WrapText = $true
$ws.Range("D4").WrapText = $true

$ws.Range("B1").WrapText = $true
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("B1").VerticalAlignment = -4108

$ws.Range("B2:B6").WrapText = $true
$ws.Range("B2:B6").HorizontalAlignment = -4108
$ws.Range("B2:B6").VerticalAlignment = -4108

# Blank formatted row 7 (A7/B7 styled like the columns above, no text)
$ws.Range("A7").WrapText = $true
$ws.Range("B7").WrapText = $true
$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("B7").VerticalAlignment = -4108

# Row heights to match the target layout
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 30

$ws.Range("D6").Select()

# Make the new sheet the active tab (clears tabSelected on the old last sheet)
$ws.Activate()

Write-Host "Added sheet 8_ with $($ws.UsedRange.Rows.Count) rows"
